$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.539.44"
$ws.Range("E2").Value = "  -0.42%  "

# Row 3
$ws.Range("D3").Value = "2.298.96"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.74%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.17%  "

# Row 9
$ws.Range("D9").Value = "2.297.25"
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("E10").Value = "  -1.32%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.95%  "

# Row 12
$ws.Range("E12").Value = "  +0.60%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.334"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.44%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.33%  "

# Row 15
$ws.Range("D15").Value = "2.707.91"
$ws.Range("E15").Value = "  +0.12%  "

# Row 16
$ws.Range("D16").Value = "58.473.63"
$ws.Range("E16").Value = "  -0.43%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.69%  "

# Row 18
$ws.Range("D18").Value = "2.291.34"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "316.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "

# Row 22
$ws.Range("E22").Value = "  +2.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "

# Row 25
$ws.Range("E25").Value = "  -1.47%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.37%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.24%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.95%  "

# Row 30
$ws.Range("E30").Value = "  -2.14%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0727"
$ws.Range("E31").Value = "  -0.45%  "

# Row 32
$ws.Range("E32").Value = "  +1.98%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.386"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "

# Row 35
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "

# Row 37
$ws.Range("E37").Value = "  -0.12%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.18%  "

# Row 39
$ws.Range("E39").Value = "  +0.41%  "

# Row 40
$ws.Range("E40").Value = "  -0.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "292.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0951"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0498"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.88%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.555"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.58%  "

# Row 48
$ws.Range("E48").Value = "  -1.98%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "

# Row 50
$ws.Range("B50").Value = "ZEEBU"
$ws.Range("C50").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.09%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
